# sl added in duplicate rows
#
# The 'Stock Report' sheet had a batch of rows (one per damage remark) whose
# cells had drifted one-off: the remark text was sitting in column X on its own
# row instead of being folded into the row above it. This re-does those rows as
# a plain [A, B] pair: A holds a single blank placeholder and B holds the raw
# serialized record (category / sub-category / remark) for that duplicate line,
# clearing out the old C:AB cells that used to carry the (wrong) per-column data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# row 7
$ws.Range("A7:AB7").Clear()
$ws.Range("A7").Value = ' '
$ws.Range("B7").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "FLOOR BOARD-(FLOOR BOARD)", "F/BOARD DIRTY BY BADLY SCRAP DUST & SCRATCHED.", "", "", ""]'

# row 11
$ws.Range("A11:AB11").Clear()
$ws.Range("A11").Value = ' '
$ws.Range("B11").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST & ODOUR ", "", "", ""]'

# row 13
$ws.Range("A13:AB13").Clear()
$ws.Range("A13").Value = ' '
$ws.Range("B13").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 15
$ws.Range("A15:AB15").Clear()
$ws.Range("A15").Value = ' '
$ws.Range("B15").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 17
$ws.Range("A17:AB17").Clear()
$ws.Range("A17").Value = ' '
$ws.Range("B17").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 19
$ws.Range("A19:AB19").Clear()
$ws.Range("A19").Value = ' '
$ws.Range("B19").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 21
$ws.Range("A21:AB21").Clear()
$ws.Range("A21").Value = ' '
$ws.Range("B21").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 24
$ws.Range("A24:AB24").Clear()
$ws.Range("A24").Value = ' '
$ws.Range("B24").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "FLOOR BOARD SLIGHTLY UP WARD.", "", "", ""]'

# row 25
$ws.Range("A25:AB25").Clear()
$ws.Range("A25").Value = ' '
$ws.Range("B25").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "FLOOR BOARD DIRTY BY DUST.", "", "", ""]'

# row 28
$ws.Range("A28:AB28").Clear()
$ws.Range("A28").Value = ' '
$ws.Range("B28").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 29
$ws.Range("A29:AB29").Clear()
$ws.Range("A29").Value = ' '
$ws.Range("B29").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B STEEL PLATE FITTING 12''X08'' & LOOSED .", "", "", ""]'

# row 31
$ws.Range("A31:AB31").Clear()
$ws.Range("A31").Value = ' '
$ws.Range("B31").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY MUD,DROP OIL SPOT,TYER MARK,SAND ,SATPLE FIBER DUST .", "", "", ""]'

# row 33
$ws.Range("A33:AB33").Clear()
$ws.Range("A33").Value = ' '
$ws.Range("B33").Value = '[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, "FLOORS-(F)", "Threshold plate-(Threshold plate)", "F/B DIRTY BY BADLY MUD+SAND,TYER MARK STAPLE FIBER DUST & ODOUR", "", "", ""]'

# Column B (now holding the long serialized-array text) needs to be much wider,
# and column X's width shrinks now that it no longer carries the long remark
# text directly. (ColumnWidth is quantized by Excel to 1/7-character pixel
# steps, so these land on the nearest representable width.)
$ws.Columns.Item(2).ColumnWidth = 251.0
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143

